# Recitation6.pptx edit: split/merge the "int i" token between the two
# "add()" code slides (Pseudocode 1 -> slide 25, C code -> slide 26).
#
# Slide 25 "Pseudocode 1" box: the line "<tab>int i = 0" loses its "int "
#   token; "<tab>" and "i" become their own runs so "i" can be flagged as
#   a (false-positive) spell-check error independently of the tab.
#
# Slide 26 "C code" box: the line "<tab>i = 0;" gains the "int " token
#   back onto the leading tab run, right before the "i" run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 25, "Pseudocode 1" content placeholder (shape id 4): paragraph 3
#   currently   Run1 = "<tab>int i"   Run2 = " = 0"
#   becomes     Run1 = "<tab>"        Run2 = "i"       Run3 = " = 0"
# ---------------------------------------------------------------------
$slide25 = $p.Slides.Item(25)
$shape25 = $slide25.Shapes.Item(1)
$para25 = $shape25.TextFrame.TextRange.Paragraphs(3)

# Prepend a new run containing just the tab; this pushes the existing
# "<tab>int i" / " = 0" runs one slot to the right (Run2 / Run3).
$para25.InsertBefore([char]9) | Out-Null

# Run2 used to read "<tab>int i" -- trim it down to just "i" so the
# leading tab (Run1) and the word "int " disappear from it, leaving the
# already-separate " = 0" run (now Run3) untouched.
$para25.Runs(2).Text = "i"

# ---------------------------------------------------------------------
# Slide 26, "C code" content placeholder (shape id 5): paragraph 4
#   currently   Run1 = "<tab>"   Run2 = "i"   Run3 = " = 0;"
#   becomes     Run1 = "<tab>int "   Run2 = "i"   Run3 = " = 0;"
# ---------------------------------------------------------------------
$slide26 = $p.Slides.Item(26)
$shape26 = $slide26.Shapes.Item(2)
$para26 = $shape26.TextFrame.TextRange.Paragraphs(4)

# Re-add "int " onto the leading tab run; Run2 ("i") and Run3 (" = 0;")
# are left exactly as they were.
$para26.Runs(1).Text = [char]9 + "int "
